$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Bad Drivers" table: D3 value correction ---
$ws.Range("D3").Value = 97.7

# --- Column width adjustments (ColumnWidth setter adds a 5/6 char offset vs. stored width) ---
$off = 5/6
$ws.Columns.Item(2).ColumnWidth = 14 - $off
$ws.Columns.Item(5).ColumnWidth = 14 - $off
$ws.Columns.Item(6).ColumnWidth = 11 - $off
$ws.Columns.Item(7).ColumnWidth = 68 - $off
$ws.Columns.Item(8).ColumnWidth = 12 - $off
$ws.Columns.Item(9).ColumnWidth = 30 - $off
$ws.Columns.Item(10).ColumnWidth = 16 - $off

# --- Rebuild the "Good Drivers" header row (row 11) with new columns & plain formatting ---
$ws.Range("A11:J11").ClearFormats()
$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- Rebuild the "Good Drivers" data row (row 12) with new columns & plain formatting ---
$ws.Range("B12:J12").ClearFormats()
$ws.Range("A12").Value = "Killer(R) Wi-Fi 6 AX1650s 160MHz Wireless Network Adapter (201D2W) - 22.250.0.4"
$ws.Range("B12").Value = 58781
$ws.Range("C12").Value = 52
$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 96
$ws.Range("F12").Value = 58842
$ws.Range("G12").Value = "killer(r) wi-fi 6 ax1650s 160mhz wireless network adapter (201d2w)"
$ws.Range("H12").Value = "22.250.0.4"
$ws.Range("I12").Value = 99.90000000000001

# J12 looks like a date ("2023-07-25"); a literal Value assignment gets auto-converted to a
# date serial by Excel's type inference. Route it through a text formula + paste-as-values so
# it lands as a plain string cell with no leftover number-format / style baggage.
$ws.Range("J12").Formula = '="2023-07-25"'
$ws.Range("J12").Copy()
$ws.Range("J12").PasteSpecial(-4163)
